# Updated cryptos list with refreshed prices and volume(1h) percentages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.534.65"
$ws.Range("E2").Value = "  +4.52%  "
$ws.Range("D3").Value = "3.249.73"
$ws.Range("E3").Value = "  +3.60%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.22"
$ws.Range("E5").Value = "  +2.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.36"
$ws.Range("E6").Value = "  +6.75%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  -3.50%  "
$ws.Range("D9").Value = "3.248.34"
$ws.Range("E9").Value = "  +3.46%  "
$ws.Range("E10").Value = "  +5.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.78"
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("E12").Value = "  +5.66%  "
$ws.Range("D13").Value = "3.820.87"
$ws.Range("E13").Value = "  +3.63%  "
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.19"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").Value = "67.510.73"
$ws.Range("E16").Value = "  +4.55%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000167"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("D18").Value = "3.258.27"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.81"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("E20").Value = "  +5.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "377.82"
$ws.Range("E21").Value = "  +6.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.61"
$ws.Range("E22").Value = "  +5.69%  "
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.27"
$ws.Range("E24").Value = "  +4.34%  "
$ws.Range("E25").Value = "  +2.66%  "
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.61"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.73"
$ws.Range("E30").Value = "  +6.99%  "
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.64"
$ws.Range("E32").Value = "  +3.44%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  +6.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.89"
$ws.Range("E35").Value = "  +4.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "162.50"
$ws.Range("E36").Value = "  +6.17%  "
$ws.Range("E37").Value = "  +4.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.853"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +7.38%  "
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.81"
$ws.Range("E40").Value = "  +3.22%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.79"
$ws.Range("E41").Value = "  +12.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.63"
$ws.Range("E42").Value = "  +4.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "363.26"
$ws.Range("E43").Value = "  +14.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.48"
$ws.Range("E44").Value = "  +7.76%  "
$ws.Range("D45").Value = "2.741.92"
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.46"
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.63"
$ws.Range("E47").Value = "  +4.32%  "
$ws.Range("E48").Value = "  +3.45%  "
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("E50").Value = "  +7.04%  "
$ws.Range("E51").Value = "  +0.77%  "
